$wb = $excel.ActiveWorkbook

# Update exchange rates on the "Exchange Rate" sheet
$rates = $wb.Worksheets.Item("Exchange Rate")
$rates.Activate()
$rates.Range("B3").Value = 1.79399
$rates.Range("B4").Value = 1.2490699999999999
$rates.Range("B6").Value = 1.6692199999999999
$rates.Range("B6").Select()
$rates.PageSetup.PaperSize = 9
$rates.PageSetup.Orientation = 1

$income = $wb.Worksheets.Item("Income")
$income.Activate()
$income.Range("D4").Select()

$costs = $wb.Worksheets.Item("Costs")
$costs.Activate()
$costs.Range("D4").Select()

$pl = $wb.Worksheets.Item("Profit and Loss")
$pl.Activate()
$pl.Range("C3").Select()

# Restore the originally active sheet/tab
$costs.Activate()

$wb.Save()
